# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price record (2022-02-09) is inserted for
# Comercializadora del Agro de Limari - Poroto granado, right before the
# existing row 32. This pushes the former rows 32..71 down to 33..72,
# growing the sheet from A1:R71 to A1:R72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 32, shifting rows 32-71 down to 33-72.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Cells.Item(32, 1).Value = 2
$ws.Cells.Item(32, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44601
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = 100112030
$ws.Cells.Item(32, 7).Value = "Poroto granado"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 500
$ws.Cells.Item(32, 11).Value = 23000
$ws.Cells.Item(32, 12).Value = 24000
$ws.Cells.Item(32, 13).Value = 23500
$ws.Cells.Item(32, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 940
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
